# Refresh the cryptos list (coin price / 1h-volume columns) as produced by the
# GitHub Actions scraper run on Thu Aug 3 08:31:25 UTC 2023.
#
# Columns: A=# (untouched), B=Coin, C=Link, D=Price, E=Volume(1h)
# Most rows just get fresh D (Price) / E (Volume) figures. A few coins
# (rows 42-47) swapped rank order, so their B/C/D/E are rewritten together.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Per-row updates: only the columns that actually changed are listed.
$updates = @(
    @{ Row = 2;  D = "29.011.94";     E = "  -1.49%  " },
    @{ Row = 3;  D = "1.828.48";      E = "  -1.48%  " },
    @{ Row = 4;  D = "1.000";         E = "  +0.01%  " },
    @{ Row = 5;  D = "238.55";        E = "  -2.61%  " },
    @{ Row = 6;  D = "0.6593";        E = "  -4.83%  " },
    @{ Row = 7;  D = "1.001";         E = "  -0.03%  " },
    @{ Row = 8;  D = "0.2942";        E = "  -3.85%  " },
    @{ Row = 9;  D = "0.07318";       E = "  -4.58%  " },
    @{ Row = 10; D = "22.63";         E = "  -3.92%  " },
    @{ Row = 11;                      E = "  -1.64%  " },
    @{ Row = 12; D = "1.832.89";      E = "  -1.24%  " },
    @{ Row = 13; D = "4.997";         E = "  -2.67%  " },
    @{ Row = 14; D = "0.6703";        E = "  -2.89%  " },
    @{ Row = 15; D = "86.01";         E = "  -4.87%  " },
    @{ Row = 16; D = "6.104";         E = "  -4.99%  " },
    @{ Row = 17; D = "29.027.00";     E = "  -1.02%  " },
    @{ Row = 18; D = "0.000008179";   E = "  -1.14%  " },
    @{ Row = 19; D = "227.10";        E = "  -4.45%  " },
    @{ Row = 20; D = "12.41";         E = "  -2.37%  " },
    @{ Row = 21; D = "1.000";         E = "  -0.02%  " },
    @{ Row = 22; D = "7.245";         E = "  -5.05%  " },
    @{ Row = 23; D = "1.001";         E = "  +0.04%  " },
    @{ Row = 24; D = "160.18";        E = "  +0.38%  " },
    @{ Row = 25; D = "0.1417";        E = "  -4.86%  " },
    @{ Row = 26; D = "8.632";         E = "  -2.94%  " },
    @{ Row = 27; D = "17.91";         E = "  -1.99%  " },
    @{ Row = 28; D = "1.498";         E = "  -2.42%  " },
    @{ Row = 29; D = "4.218";         E = "  -0.80%  " },
    @{ Row = 30; D = "4.103";         E = "  -1.23%  " },
    @{ Row = 31; D = "1.198";         E = "  +0.47%  " },
    @{ Row = 32; D = "0.05324";       E = "  +4.16%  " },
    @{ Row = 33; D = "0.7459";        E = "  -2.76%  " },
    @{ Row = 34; D = "1.840";         E = "  -2.38%  " },
    @{ Row = 35; D = "1.124";         E = "  -2.27%  " },
    @{ Row = 36; D = "2.683";         E = "  -0.03%  " },
    @{ Row = 37; D = "1.294.52";      E = "  -2.58%  " },
    @{ Row = 38;                      E = "  -3.12%  " },
    @{ Row = 39; D = "2.700";         E = "  -0.63%  " },
    @{ Row = 40; D = "0.9205";        E = "  -6.08%  " },
    @{ Row = 41; D = "6.014";         E = "  +3.10%  " },

    # Rank shuffle: Quant/PaxDollar swapped; XinFinNetwork, RocketPoolETH,
    # BabyDogeCoin and Mantle rotated by one slot.
    @{ Row = 42; B = "PaxDollar";      C = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp";         D = "0.9994";         E = "  -0.11%  " },
    @{ Row = 43; B = "Quant";          C = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt";          D = "103.61";         E = "  -2.23%  " },
    @{ Row = 44; B = "RocketPoolETH";  C = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth";     D = "1.979.71";       E = "  -1.30%  " },
    @{ Row = 45; B = "BabyDogeCoin";   C = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge";  D = "0.00000000123";  E = "  -2.33%  " },
    @{ Row = 46; B = "Mantle";         C = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt";             D = "0.5176";         E = "  -0.83%  " },
    @{ Row = 47; B = "XinFinNetwork";  C = "https://coinranking.com/coin/77jGXSqWJ1ofG+xinfinnetwork-xdc";  D = "0.07734";        E = "  +16.98%  " },

    @{ Row = 48; D = "63.20";         E = "  +0.39%  " },
    @{ Row = 49; D = "1.744";         E = "  -1.68%  " },
    @{ Row = 50; D = "9.259" },
    @{ Row = 51;                      E = "  -0.25%  " }
)

foreach ($u in $updates) {
    $row = $u.Row
    foreach ($col in @("B", "C", "D", "E")) {
        if ($u.ContainsKey($col)) {
            $cell = $ws.Range("$col$row")
            # Force text so things like "1.000" / "29.011.94" / "  -1.49%  "
            # are stored verbatim instead of being re-interpreted as numbers.
            $cell.NumberFormat = "@"
            $cell.Value = $u[$col]
        }
    }
}
